$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$newRows = @(
    @{
        Row = 112
        Idx = 111
        B = "paraguay"
        C = "primera-division"
        D = "2023"
        E = 45236.91666666666
        F = "Libertad Asuncion"
        G = 1
        H = "Cerro Porteno"
        I = 1
        J = 1.96
        K = "04/11/2023 00:42"
        L = 2.03
        M = "06/11/2023 21:53"
        N = 3.63
        O = "04/11/2023 00:42"
        P = 3.42
        Q = "06/11/2023 21:54"
        R = 3.86
        S = "04/11/2023 00:42"
        T = 4.02
        U = "06/11/2023 21:53"
        V = "https://www.betexplorer.com/football/paraguay/primera-division/libertad-asuncion-cerro-porteno/EZ4G0FS6/"
    },
    @{
        Row = 113
        Idx = 112
        B = "paraguay"
        C = "primera-division"
        D = "2023"
        E = 45237.02083333334
        F = "Nacional Asuncion"
        G = 5
        H = "Resistencia"
        I = 0
        J = 1.45
        K = "02/11/2023 22:12"
        L = 1.39
        M = "06/11/2023 23:59"
        N = 4.44
        O = "02/11/2023 22:12"
        P = 4.83
        Q = "06/11/2023 23:59"
        R = 6.37
        S = "02/11/2023 22:12"
        T = 8.619999999999999
        U = "07/11/2023 00:21"
        V = "https://www.betexplorer.com/football/paraguay/primera-division/nacional-asuncion-resistencia/foMDaZsD/"
    }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Column A (Indice) - numeric, bold/bordered style matching the row above
    $ws.Cells.Item($row, 1).Value = $r.Idx
    $ws.Range("A$prevRow").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial($xlPasteFormats) | Out-Null

    # Columns B-D - plain text
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).Style = "Normal"

    # Column E (data_partida) - numeric date/time, matching style of row above
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Range("E$prevRow").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial($xlPasteFormats) | Out-Null

    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}
